$p = $ppt.ActivePresentation
$master = $p.Designs.Item(1).SlideMaster
$layout = $master.CustomLayouts.Item(1)
$tcs = $layout.ThemeColorScheme
Write-Output "Name: $($tcs.Name)"
Write-Output "Count: $($tcs.Count)"
